$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 62, pushing existing rows 62:84 down to 63:85
$ws.Rows.Item(62).Insert()

# Populate the newly inserted row 62 with the new record
$ws.Cells.Item(62, 1).Value = 6
$ws.Cells.Item(62, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(62, 3).Value = "Metropolitana"
$ws.Cells.Item(62, 4).Value = 44809
$ws.Cells.Item(62, 5).Value = 13
$ws.Cells.Item(62, 6).Value = "Fruta"
$ws.Cells.Item(62, 7).Value = 100108
$ws.Cells.Item(62, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(62, 9).Value = 100108007
$ws.Cells.Item(62, 10).Value = "Coco"
$ws.Cells.Item(62, 11).Value = "Sin especificar"
$ws.Cells.Item(62, 12).Value = "Primera"
$ws.Cells.Item(62, 13).Value = 50
$ws.Cells.Item(62, 14).Value = 27000
$ws.Cells.Item(62, 15).Value = 28000
$ws.Cells.Item(62, 16).Value = 27500
$ws.Cells.Item(62, 17).Value = "`$/malla 20 unidades"
$ws.Cells.Item(62, 18).Value = "Perú"
$ws.Cells.Item(62, 19).Value = 1375
$ws.Cells.Item(62, 20).Value = 20
